# Apply refreshed cryptocurrency market data (prices / 1h volume %) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.382.67"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "2.876.44"
$ws.Range("E3").Value = "  +7.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "195.97"
$ws.Range("E5").Value = "  +4.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "598.26"
$ws.Range("E6").Value = "  +1.73%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.554"
$ws.Range("E8").Value = "  +3.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.193"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("D10").Value = "2.873.29"
$ws.Range("E10").Value = "  +7.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.392"
$ws.Range("E11").Value = "  +9.86%  "
$ws.Range("E12").Value = "  -1.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.91"
$ws.Range("E13").Value = "  +4.10%  "
$ws.Range("D14").Value = "3.397.12"
$ws.Range("E14").Value = "  +7.57%  "
$ws.Range("D15").Value = "76.151.96"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.55"
$ws.Range("E16").Value = "  +3.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000189"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "2.877.57"
$ws.Range("E18").Value = "  +7.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.03"
$ws.Range("E19").Value = "  -2.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.55"
$ws.Range("E20").Value = "  +5.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.35"
$ws.Range("E21").Value = "  +2.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.34"
$ws.Range("E22").Value = "  +2.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.13"
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.80"
$ws.Range("E24").Value = "  +2.78%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "3.028.82"
$ws.Range("E26").Value = "  +7.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.22"
$ws.Range("E27").Value = "  +1.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.75"
$ws.Range("E28").Value = "  +4.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000105"
$ws.Range("E29").Value = "  +10.38%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.41"
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "510.78"
$ws.Range("E32").Value = "  -1.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.71"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.81"
$ws.Range("E34").Value = "  +3.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "167.54"
$ws.Range("E36").Value = "  +2.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.02"
$ws.Range("E37").Value = "  +4.39%  "
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.53"
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "184.67"
$ws.Range("E40").Value = "  +8.43%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.345"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.06"
$ws.Range("E43").Value = "  +1.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.68"
$ws.Range("E44").Value = "  -1.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0923"
$ws.Range("E45").Value = "  +9.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.23"
$ws.Range("E46").Value = "  +2.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.23"
$ws.Range("E47").Value = "  +2.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.36"
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.688"
$ws.Range("E49").Value = "  +16.24%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.580"
$ws.Range("E50").Value = "  +8.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.75"
$ws.Range("E51").Value = "  +2.84%  "
